$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Testdata")

$ws2.Range("A6").Value = "EleType1"
$ws2.Range("A7").Value = "EleType2"
$ws2.Range("B6").Value = "JSElement"
$ws2.Range("B7").Value = "JSElement"

$ws2.Range("B3:B5").Borders.LineStyle = 1
$ws2.Range("A6:B7").Borders.LineStyle = 1
